# Auto-generated edit script: apply cached-value updates to Seraph_Profits sheets
$wb = $excel.ActiveWorkbook

# --- Sheet ALC (51 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2119.8
$ws.Range("I9").Value = 1799.6666
$ws.Range("J9").Value = 2600
$ws.Range("K9").Value = 1799.6666
$ws.Range("L9").Value = 2600
$ws.Range("N9").Value = -2938
$ws.Range("M9").Value = -1630.6666
$ws.Range("H32").Value = 874.75
$ws.Range("I32").Value = 799.5
$ws.Range("J32").Value = 950
$ws.Range("K32").Value = 799.5
$ws.Range("L32").Value = 950
$ws.Range("M32").Value = -473.5
$ws.Range("N32").Value = -1602
$ws.Range("H80").Value = 19287.375
$ws.Range("I80").Value = 299.83334
$ws.Range("J80").Value = 76250
$ws.Range("K80").Value = 899.5000200000001
$ws.Range("L80").Value = 228750
$ws.Range("M80").Value = 98.49997999999994
$ws.Range("N80").Value = -230746
$ws.Range("H83").Value = 19287.375
$ws.Range("I83").Value = 299.83334
$ws.Range("J83").Value = 76250
$ws.Range("K83").Value = 2698.50006
$ws.Range("L83").Value = 686250
$ws.Range("M83").Value = 2293.49994
$ws.Range("N83").Value = -696234
$ws.Range("H86").Value = 3706.1538
$ws.Range("I86").Value = 2199.5
$ws.Range("K86").Value = 2199.5
$ws.Range("M86").Value = -1076.5
$ws.Range("H89").Value = 3706.1538
$ws.Range("I89").Value = 2199.5
$ws.Range("K89").Value = 10997.5
$ws.Range("M89").Value = -5381.5
$ws.Range("H112").Value = 3378.1667
$ws.Range("J112").Value = 3556.182
$ws.Range("L112").Value = 10668.546
$ws.Range("N112").Value = -12884.546
$ws.Range("H116").Value = 9000
$ws.Range("J116").Value = 9000
$ws.Range("L116").Value = 9000
$ws.Range("N116").Value = -15884
$ws.Range("H138").Value = 7393.619
$ws.Range("I138").Value = 9134.909
$ws.Range("J138").Value = 5478.2
$ws.Range("K138").Value = 27404.727
$ws.Range("L138").Value = 16434.6
$ws.Range("M138").Value = -22264.727
$ws.Range("N138").Value = -26714.6

# --- Sheet ARM (31 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2600
$ws.Range("J2").Value = 1950
$ws.Range("L2").Value = 1950
$ws.Range("N2").Value = -2176
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H110").Value = 4758.357
$ws.Range("I110").Value = 5051.5
$ws.Range("J110").Value = 2999.5
$ws.Range("K110").Value = 5051.5
$ws.Range("L110").Value = 2999.5
$ws.Range("M110").Value = -3006.5
$ws.Range("N110").Value = -7089.5
$ws.Range("H116").Value = 2600
$ws.Range("J116").Value = 1950
$ws.Range("L116").Value = 1950
$ws.Range("N116").Value = -6538
$ws.Range("H119").Value = 40000
$ws.Range("J119").Value = 40000
$ws.Range("L119").Value = 40000
$ws.Range("N119").Value = -49676
$ws.Range("H122").Value = 2988.4614
$ws.Range("I122").Value = 2176.353
$ws.Range("K122").Value = 6529.059
$ws.Range("M122").Value = -4079.059
$ws.Range("H132").Value = 9505.076999999999
$ws.Range("I132").Value = 7855.6665
$ws.Range("K132").Value = 23566.9995
$ws.Range("M132").Value = -21036.9995

# --- Sheet BSM (22 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2600
$ws.Range("J3").Value = 1950
$ws.Range("L3").Value = 1950
$ws.Range("N3").Value = -2178
$ws.Range("H86").Value = 2796.2727
$ws.Range("I86").Value = 2076
$ws.Range("J86").Value = 4056.75
$ws.Range("K86").Value = 2076
$ws.Range("L86").Value = 4056.75
$ws.Range("M86").Value = -953
$ws.Range("N86").Value = -6302.75
$ws.Range("H89").Value = 2796.2727
$ws.Range("I89").Value = 2076
$ws.Range("J89").Value = 4056.75
$ws.Range("K89").Value = 10380
$ws.Range("L89").Value = 20283.75
$ws.Range("M89").Value = -4764
$ws.Range("N89").Value = -31515.75
$ws.Range("H134").Value = 2072.3125
$ws.Range("I134").Value = 1582.9286
$ws.Range("K134").Value = 4748.7858
$ws.Range("M134").Value = -2213.7858

# --- Sheet CRP (8 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13933.441
$ws.Range("J99").Value = 17779.455
$ws.Range("L99").Value = 17779.455
$ws.Range("N99").Value = -20775.455
$ws.Range("H126").Value = 13933.441
$ws.Range("J126").Value = 17779.455
$ws.Range("L126").Value = 53338.36500000001
$ws.Range("N126").Value = -58278.36500000001

# --- Sheet CUL (12 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 2235.5557
$ws.Range("I14").Value = 2235.5557
$ws.Range("K14").Value = 6706.6671
$ws.Range("M14").Value = -6533.6671
$ws.Range("H121").Value = 1836.5714
$ws.Range("J121").Value = 4999
$ws.Range("L121").Value = 14997
$ws.Range("N121").Value = -17617
$ws.Range("H134").Value = 5066
$ws.Range("I134").Value = 1099.5
$ws.Range("K134").Value = 3298.5
$ws.Range("M134").Value = 1771.5

# --- Sheet GSM (7 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2281.1667
$ws.Range("I132").Value = 959.8
$ws.Range("J132").Value = 8888
$ws.Range("K132").Value = 2879.4
$ws.Range("L132").Value = 26664
$ws.Range("M132").Value = -349.3999999999996
$ws.Range("N132").Value = -31724

# --- Sheet LTW (36 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1090.4
$ws.Range("J22").Value = 1404.3334
$ws.Range("L22").Value = 1404.3334
$ws.Range("N22").Value = -1994.3334
$ws.Range("H27").Value = 1090.4
$ws.Range("J27").Value = 1404.3334
$ws.Range("L27").Value = 1404.3334
$ws.Range("N27").Value = -1618.3334
$ws.Range("H40").Value = 4634.9287
$ws.Range("I40").Value = 4582.5
$ws.Range("K40").Value = 4582.5
$ws.Range("M40").Value = -4446.5
$ws.Range("H61").Value = 7256.8
$ws.Range("I61").Value = 10696.75
$ws.Range("K61").Value = 10696.75
$ws.Range("M61").Value = -10494.75
$ws.Range("H64").Value = 34944
$ws.Range("J64").Value = 34944
$ws.Range("L64").Value = 34944
$ws.Range("N64").Value = -35394
$ws.Range("H67").Value = 34944
$ws.Range("J67").Value = 34944
$ws.Range("L67").Value = 34944
$ws.Range("N67").Value = -36504
$ws.Range("H68").Value = 3387.8572
$ws.Range("I68").Value = 3405
$ws.Range("K68").Value = 3405
$ws.Range("M68").Value = -2656
$ws.Range("H71").Value = 3387.8572
$ws.Range("I71").Value = 3405
$ws.Range("K71").Value = 17025
$ws.Range("M71").Value = -13281
$ws.Range("H113").Value = 7256.8
$ws.Range("I113").Value = 10696.75
$ws.Range("K113").Value = 10696.75
$ws.Range("M113").Value = -8526.75

# --- Sheet WVR (8 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 78599.8
$ws.Range("J46").Value = 78599.8
$ws.Range("L46").Value = 78599.8
$ws.Range("N46").Value = -79061.8
$ws.Range("H134").Value = 78599.8
$ws.Range("J134").Value = 78599.8
$ws.Range("L134").Value = 235799.4
$ws.Range("N134").Value = -240869.4

